$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# Add two new data rows (21 and 22) to the "Bluff Lake Sondes" table.
#
# Row 21 mirrors the existing Cypress Boardwalk / 30min entries (e.g. row 7),
# and row 22 mirrors the bottom-bordered Visitor Center Baro entry that closes
# that block (e.g. row 11). Formatting is copied from those analogous rows
# first so the new cells reuse the workbook's existing styles/number formats
# instead of creating brand-new ones.
# ---------------------------------------------------------------------------

$ws.Range("A7:J7").Copy() | Out-Null
$ws.Range("A21:J21").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

# Column E was blank on row 7 (so it had no explicit style); row 10 carries
# the "time value, no border" style that column E needs here.
$ws.Range("E10").Copy() | Out-Null
$ws.Range("E21").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$ws.Range("A11:G11").Copy() | Out-Null
$ws.Range("A22:G22").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

# --- Row 21 values: Cypress Boardwalk sonde, checked 2019-05-07 -----------
$ws.Cells.Item(21, 1).Value = 43592
$ws.Cells.Item(21, 2).Value = 20089579
$ws.Cells.Item(21, 3).Value = "Cypress Boardwalk"
$ws.Cells.Item(21, 4).Value = "Good"
$ws.Cells.Item(21, 5).Value = 0.58333333333333337
$ws.Cells.Item(21, 6).Value = 0.625
$ws.Cells.Item(21, 7).Value = "30min"
$ws.Cells.Item(21, 8).Value = 0.2
$ws.Cells.Item(21, 9).Value = 1.5
$ws.Cells.Item(21, 10).Value = 1.1200000000000001

# --- Row 22 values: Visitor Center Baro sonde, checked 2019-08-07 ---------
$ws.Cells.Item(22, 1).Value = 43684
$ws.Cells.Item(22, 2).Value = 10868630
$ws.Cells.Item(22, 3).Value = "Visitor Center Baro"
$ws.Cells.Item(22, 4).Value = "Good"
$ws.Cells.Item(22, 5).Value = 0.66666666666666663
$ws.Cells.Item(22, 6).Value = 0.70833333333333337
$ws.Cells.Item(22, 7).Value = "30min"

# ---------------------------------------------------------------------------
# Update the view state to match: scrolled down one row, selection on F21.
# ---------------------------------------------------------------------------
$excel.ActiveWindow.ScrollRow = 2
$excel.ActiveWindow.ScrollColumn = 3
$ws.Range("F21").Select()
